# Re-curate the dimension/measure metadata for columns B (tipo-de-vehiculo-orden),
# C (tipo-de-vehiculo) and E (municipio-nombre):
#  - B/C switch from dimensions to measures (iaest-dimension:* -> iaest-measure:*,
#    "dim" -> "medida", "skos:Concept" -> "xsd:int") and lose their per-column
#    mapping file (row 5, B5/C5 cleared).
#  - E switches from a measure to a dimension (iaest-measure:municipio-nombre ->
#    sdmx-dimension:refArea, "medida" -> "dim", and the URI column now points to
#    a new "URI-Municipio" value instead of "xsd:int").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sdmx/iaest dimension-or-measure identifiers
$ws.Range("B2").Value = "iaest-measure:tipo-de-vehiculo-orden"
$ws.Range("C2").Value = "iaest-measure:tipo-de-vehiculo"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: dim / medida flag
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("E3").Value = "dim"

# Row 4: data type / URI mapping
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Municipio"

# Row 5: per-column external mapping workbook references — B5/C5 no longer apply,
# remove the cells outright (not just their value) to match the curated layout.
$ws.Range("B5:C5").Clear()
